$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.522.99"
$ws.Range("E2").Value = "  +6.95%  "
$ws.Range("D3").Value = "1.728.38"
$ws.Range("E3").Value = "  +3.88%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.88"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3742"
$ws.Range("E7").Value = "  +2.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.54"
$ws.Range("E8").Value = "  +2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3381"
$ws.Range("E9").Value = "  +4.03%  "
$ws.Range("E10").Value = "  +4.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07474"
$ws.Range("E11").Value = "  +5.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.414"
$ws.Range("E13").Value = "  +5.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.16"
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.063"
$ws.Range("E15").Value = "  +7.22%  "
$ws.Range("D16").Value = "1.726.72"
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001076"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06652"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.05"
$ws.Range("E19").Value = "  +4.56%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.57"
$ws.Range("E21").Value = "  +4.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.166"
$ws.Range("E22").Value = "  +4.07%  "
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "26.518.38"
$ws.Range("E24").Value = "  +6.99%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.421"
$ws.Range("E26").Value = "  +22.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.403"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.24"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.45"
$ws.Range("E29").Value = "  +4.40%  "
$ws.Range("D30").Value = "1.914.39"
$ws.Range("E30").Value = "  +3.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.36"
$ws.Range("E31").Value = "  +4.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.096"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.997"
$ws.Range("E33").Value = "  +5.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08652"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.697"
$ws.Range("E35").Value = "  +3.71%  "
$ws.Range("E36").Value = "  +5.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.379"
$ws.Range("E37").Value = "  +4.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02338"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2157"
$ws.Range("E39").Value = "  +3.35%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2157"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.426"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.219"
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6216"
$ws.Range("E43").Value = "  +5.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.28"
$ws.Range("E44").Value = "  +7.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.845"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6036"
$ws.Range("E47").Value = "  +6.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.98"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("E49").Value = "  +5.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07184"
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.10"
$ws.Range("E51").Value = "  +2.79%  "
